$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

$find.Execute(
    "môžete pozorovať súhvezdie ozvezdje Bik 2022: 16.-25. januarja",
    $true,
    $false,
    $false,
    $false,
    $false,
    $true,
    1,
    $false,
    "2022: Datumi kampanje za opazovanje ozvezdje Bik: 16.-25. januarja",
    2
)
